$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1255
$ws.Range("C2").Value = 5247.01035856574
$ws.Range("D2").Value = 4259.68905145549
$ws.Range("E2").Value = 6234.33166567599
$ws.Range("F2").Value = 53.7238678487814
$ws.Range("G2").Value = 43.3349333927245
$ws.Range("H2").Value = 64.865793613918

$ws.Range("B3").Value = 106
$ws.Range("C3").Value = 4768.32075471698
$ws.Range("D3").Value = 1246.6639642203
$ws.Range("E3").Value = 8289.97754521366
$ws.Range("F3").Value = 29.9533838592261
$ws.Range("G3").Value = 6.79666344515095
$ws.Range("H3").Value = 58.1311759345058

$ws.Range("B4").Value = 212
$ws.Range("C4").Value = 5083.29716981132
$ws.Range("D4").Value = 1925.62395561542
$ws.Range("E4").Value = 8240.97038400722
$ws.Range("F4").Value = 37.5309408963742
$ws.Range("G4").Value = 15.673427141843
$ws.Range("H4").Value = 63.5186245553872

$ws.Range("B5").Value = 2
$ws.Range("C5").Value = -8218
$ws.Range("D5").Value = -67012.5146420991
$ws.Range("E5").Value = 50576.5146420991
$ws.Range("F5").Value = -66.1499228684452
$ws.Range("G5").Value = -98.2917231587202
$ws.Range("H5").Value = 570.750603253392

$ws.Range("B6").Value = 463
$ws.Range("C6").Value = 4223.86177105832
$ws.Range("D6").Value = 2385.15949751633
$ws.Range("E6").Value = 6062.5640446003
$ws.Range("F6").Value = 36.7560429998814
$ws.Range("G6").Value = 22.5634156306051
$ws.Range("H6").Value = 52.5921515874865

$ws.Range("B7").Value = 270
$ws.Range("C7").Value = 3486.8962962963
$ws.Range("D7").Value = 1130.77731253732
$ws.Range("E7").Value = 5843.01528005527
$ws.Range("F7").Value = 28.7015248956852
$ws.Range("G7").Value = 12.1684636931031
$ws.Range("H7").Value = 47.6714752534599

$ws.Range("B8").Value = 165
$ws.Range("C8").Value = 5303.29696969697
$ws.Range("D8").Value = 1875.37084991641
$ws.Range("E8").Value = 8731.22308947753
$ws.Range("F8").Value = 38.9557843140747
$ws.Range("G8").Value = 14.9414938196058
$ws.Range("H8").Value = 67.9872894695763

$ws.Range("B9").Value = 435
$ws.Range("C9").Value = 5005.44827586207
$ws.Range("D9").Value = 3324.84368519627
$ws.Range("E9").Value = 6686.05286652787
$ws.Range("F9").Value = 39.6853044517943
$ws.Range("G9").Value = 25.2680763592285
$ws.Range("H9").Value = 55.7618257331291

$ws.Range("B10").Value = 501
$ws.Range("C10").Value = 6611.32734530938
$ws.Range("D10").Value = 4807.64893496511
$ws.Range("E10").Value = 8415.00575565366
$ws.Range("F10").Value = 61.8168318575774
$ws.Range("G10").Value = 44.7625769896183
$ws.Range("H10").Value = 80.8802220638922

$ws.Range("B11").Value = 461
$ws.Range("C11").Value = 9293.39913232104
$ws.Range("D11").Value = 7228.49207611551
$ws.Range("E11").Value = 11358.3061885266
$ws.Range("F11").Value = 77.5932530158467
$ws.Range("G11").Value = 57.3205814731995
$ws.Range("H11").Value = 100.47830500883

$ws.Range("B12").Value = 631
$ws.Range("C12").Value = 3567.27733755943
$ws.Range("D12").Value = 2278.17299178467
$ws.Range("E12").Value = 4856.38168333419
$ws.Range("F12").Value = 33.720617903291
$ws.Range("G12").Value = 21.8854097640282
$ws.Range("H12").Value = 46.7050378470745

$ws.Range("B13").Value = 289
$ws.Range("C13").Value = 8897.52249134948
$ws.Range("D13").Value = 6486.7576265608
$ws.Range("E13").Value = 11308.2873561382
$ws.Range("F13").Value = 86.3386600004225
$ws.Range("G13").Value = 59.7818640955245
$ws.Range("H13").Value = 117.309369916943

$ws.Range("C14").Value = 4325.71428571429
$ws.Range("D14").Value = 1267.34382334928
$ws.Range("E14").Value = 7384.08474807929
$ws.Range("F14").Value = 37.9133889131261
$ws.Range("G14").Value = 16.532889283864
$ws.Range("H14").Value = 63.2166074177725
